$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '47.140.48'
$ws.Range('D3').Value = '2.495.33'
$ws.Range('E3').Value = '  +1.75%  '
$ws.Range('E4').Value = '  +0.03%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '322.56'
$c.ClearFormats()
$ws.Range('E5').Value = '  +0.78%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '105.31'
$c.ClearFormats()
$ws.Range('E6').Value = '  +0.67%  '
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.ClearFormats()
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  +1.84%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '37.54'
$c.ClearFormats()
$ws.Range('E10').Value = '  +4.11%  '
$ws.Range('E11').Value = '  +0.91%  '
$ws.Range('E12').Value = '  +0.14%  '
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '18.40'
$c.ClearFormats()
$ws.Range('E13').Value = '  -1.42%  '
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '7.23'
$c.ClearFormats()
$ws.Range('E14').Value = '  +2.42%  '
$ws.Range('D15').Value = '2.885.34'
$ws.Range('E15').Value = '  +1.68%  '
$ws.Range('D16').Value = '2.492.25'
$ws.Range('E16').Value = '  +1.90%  '
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '0.845'
$c.ClearFormats()
$ws.Range('E17').Value = '  +0.51%  '
$ws.Range('D18').Value = '47.061.77'
$ws.Range('E18').Value = '  +4.05%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '12.63'
$c.ClearFormats()
$ws.Range('E19').Value = '  +1.75%  '
$ws.Range('E20').Value = '  +2.47%  '
$ws.Range('E21').Value = '  +0.29%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '70.91'
$c.ClearFormats()
$ws.Range('E22').Value = '  +2.41%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '251.21'
$c.ClearFormats()
$ws.Range('E23').Value = '  +2.77%  '
$ws.Range('E24').Value = '  +2.70%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '2.55'
$c.ClearFormats()
$ws.Range('E25').Value = '  +0.75%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '26.15'
$c.ClearFormats()
$ws.Range('E26').Value = '  +2.41%  '
$ws.Range('E27').Value = '  -0.08%  '
$ws.Range('E28').Value = '  +5.94%  '
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '2.19'
$c.ClearFormats()
$ws.Range('E29').Value = '  -0.16%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '35.30'
$c.ClearFormats()
$ws.Range('E30').Value = '  +3.95%  '
$ws.Range('E31').Value = '  +4.05%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '49.66'
$c.ClearFormats()
$ws.Range('E32').Value = '  -0.28%  '
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '19.77'
$c.ClearFormats()
$ws.Range('E33').Value = '  -3.19%  '
$ws.Range('E34').Value = '  +2.10%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '0.0783'
$c.ClearFormats()
$ws.Range('E35').Value = '  +1.98%  '
$ws.Range('E36').Value = '  +0.08%  '
$ws.Range('E37').Value = '  +0.83%  '
$ws.Range('E38').Value = '  +1.84%  '
$ws.Range('E39').Value = '  +3.31%  '
$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '0.111'
$c.ClearFormats()
$ws.Range('E40').Value = '  +0.97%  '
$ws.Range('B41').Value = 'Monero'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '121.79'
$c.ClearFormats()
$ws.Range('E41').Value = '  -2.89%  '
$ws.Range('E42').Value = '  +1.06%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '21.75'
$c.ClearFormats()
$ws.Range('E43').Value = '  +1.45%  '
$ws.Range('E44').Value = '  +0.82%  '
$ws.Range('D45').Value = '1.952.61'
$ws.Range('E45').Value = '  -0.04%  '
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '2.98'
$c.ClearFormats()
$ws.Range('E46').Value = '  -0.29%  '
$ws.Range('E47').Value = '  -0.37%  '
$ws.Range('E48').Value = '  -1.12%  '
$ws.Range('E49').Value = '  -0.13%  '
$ws.Range('E50').Value = '  +12.94%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '78.77'
$c.ClearFormats()
$ws.Range('E51').Value = '  +3.02%  '
